# Update build version string from "January 30 2026 16.19.47 EST" to
# "February 02 2026 12.49.33 EST" throughout the workbook.

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: version banner and recommended citation text.
$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Shanjiaoshu Coal Mine, China, M2110, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# "Boundaries and methane sources" sheet: build_version column (S), rows 2-8.
for ($row = 2; $row -le 8; $row++) {
    $cell = $wsData.Range("S" + $row)
    $cell.Value = "mines - January 30 (built on " + $newStamp + ")"
}
